$d = $word.ActiveDocument

$pairs = @(
    @{old="357×6=2142"; new="827×6=4962"},
    @{old="689×2=1378"; new="934×9=8406"},
    @{old="810×5=4050"; new="129×3=387"},
    @{old="407×5=2035"; new="265×4=1060"},
    @{old="987×8=7896"; new="739×8=5912"},
    @{old="623×2=1246"; new="866×8=6928"},
    @{old="590×9=5310"; new="960×9=8640"},
    @{old="549×6=3294"; new="744×6=4464"},
    @{old="752×2=1504"; new="975×8=7800"},
    @{old="886×3=2658"; new="225×3=675"},
    @{old="609×6=3654"; new="678×9=6102"},
    @{old="314×7=2198"; new="527×9=4743"},
    @{old="470×8=3760"; new="913×6=5478"},
    @{old="250×7=1750"; new="898×6=5388"},
    @{old="978×6=5868"; new="867×5=4335"},
    @{old="771×7=5397"; new="893×4=3572"},
    @{old="201×7=1407"; new="987×5=4935"},
    @{old="209×8=1672"; new="208×8=1664"},
    @{old="675×8=5400"; new="552×4=2208"},
    @{old="861×2=1722"; new="171×7=1197"},
    @{old="131×4=524";  new="781×2=1562"},
    @{old="403×2=806";  new="447×4=1788"},
    @{old="589×8=4712"; new="900×8=7200"},
    @{old="942×6=5652"; new="538×2=1076"},
    @{old="901×6=5406"; new="795×2=1590"}
)

foreach ($p in $pairs) {
    $d.Content.Find.Execute($p.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $p.new, 2)
}
